$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of row 3's cells: A3 and the previously-empty B3/C3,
# turning the "Next/Previous cell is empty" story into
# "The above row" / "is full" / "of nones".
$ws.Range("A3").Value = "The above row"
$ws.Range("B3").Value = "is full"
$ws.Range("C3").Value = "of nones"

# Move the active selection from A3 to C3.
$ws.Range("C3").Select()
